$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Potencia sheet: relabel the EU28 data source as EU27 and refresh the
# underlying generation / consumption data series with the updated figures.
# ---------------------------------------------------------------------------
$wsPotencia = $wb.Worksheets.Item("Potencia")
$wsPotencia.Range("A2").Value = "EU27 - Electricity balance (GWh)"

# ---------------------------------------------------------------------------
# About sheet: insert two new rows (File + Tab lines) after the existing
# "Tab: Electricity Balance; rows: 25&26" row, pushing the Notes section down.
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("8:9").EntireRow.Insert()
$wsAbout.Range("B8").Value = "File: Central_2018_UK_pg_det_yearly"
$wsAbout.Range("B8").HorizontalAlignment = -4131
$wsAbout.Range("B9").Value = "Tab: Electricity Balance; rows: 25&26"

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY","AZ")
$row3vals = @(185890.91902720978,187304.32674418611,180361.63186046516,185252.41023255815,187650.86953488376,192210.59360439423,186517.81279069762,189072.72418604651,190770.98255813954,180797.16976744187,185513.99718863514,177602.2230626955,181445.9325442444,183651.33326543003,174654.1504688836,178359.8051451685,181294.57396427891,181823.1999539338,180506.60896212255,178988.98644971329,178014.04951468456,177607.14549645549,177379.67694948078,175985.34795745445,175119.91302955168,174362.60281991915,173954.80174730317,173472.15459573959,173370.91484976248,173253.8111454621,172800.91296743666,172191.83430205891,172677.28785541485,172782.49615009592,172561.0161748291,172817.92214755566,173329.49017028921,174129.50937778977,175045.18253572617,176085.85861473411,177256.32448998454,178804.68854397655,180183.12943058135,181343.125363075,182456.86203022164,183536.15448119456,184813.44232076689,186034.1076974585,187024.11666272618,188215.08481829744,189297.8606713576)
$row4vals = @(2198618.1929556685,2260345.821162791,2288731.9329069769,2349583.9139534878,2407136.116976744,2434893.6473607589,2489359.4653488379,2509541.9643023256,2522423.9036046509,2389095.1697674412,2510932.53195599,2464110.9670211952,2474306.4779876219,2450910.5725532961,2403562.1852317555,2439191.9123678105,2460798.6261265981,2480405.0569103057,2476330.5099030756,2474449.2737449058,2483321.4086599997,2502085.5751254796,2525629.3410594575,2534703.6260747798,2550961.840323214,2567522.4495825293,2587976.2721975846,2605629.8216055217,2626638.3077441482,2643343.0344365947,2654101.653601754,2659696.4271210991,2679151.0539736128,2691642.1998726409,2698571.1478789761,2711212.9601077675,2727281.3944133762,2746350.4898617445,2766724.200705783,2787857.0217624293,2811739.7895229333,2839773.1096237986,2865244.8594552097,2887420.2803410199,2910516.1497898763,2932467.8190216823,2958147.5852745809,2983982.7922188127,3004663.1641574157,3029362.4441876449,3053060.6868686839)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $wsPotencia.Range($cols[$i] + "3").Value = $row3vals[$i]
    $wsPotencia.Range($cols[$i] + "4").Value = $row4vals[$i]
}
